# Generate Report for Handoff
# Updates the Overview / zh-cn / de-de sheets of the localization-status
# workbook: the markdown source file was swapped for a new UUID, and two
# new image assets (with their own zh-cn/de-de handoff rows) were added,
# plus refreshed handoff timestamps for the markdown row.

$wb = $excel.ActiveWorkbook

$mdOld   = "8177b148-a9e4-41f9-b46f-e2e11509f714"
$mdNew   = "53be77f6-6bb5-4b53-8408-d7a0e180bd9f"
$hashNew = "e1daab25dfff109b568ac1c1a2e0a116b6be8800"

$img1 = "5fabb73f-d69e-4802-bf8e-6f7fd29c6902.png"
$img2 = "6809a36d-22f3-4b4b-93e1-f17256f8ecf3.png"
$img1Handoff = "663623632995ee455910745046fff25f4d69b8dd.png"
$img2Handoff = "413a74eb561cda36a8e1d172a881db8c2c25a125.png"

$mdFile  = "$mdNew.md"
$xlfZh   = "$mdNew.$hashNew.zh-cn.xlf"
$xlfDe   = "$mdNew.$hashNew.de-de.xlf"
$depFrom = "e2e\$mdFile"

$zeroDate = "0001-01-01 00:00:00"
$zhTime   = "2016-03-08 21:16:05"
$deTime   = "2016-03-08 21:16:16"

$repoBase    = "https://github.com/OpenLocalizationTest/oltest/blob/3cddcfe616b9f1b2246cb17065695b475c0dbbfd"
$e2eBase     = "$repoBase/e2e"
$handoffZhBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/55f949617cc16286aa65d8d4dcbe6e7f569c5da6/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht"
$handoffDeBase = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/79bfb952ea3e99e4a0e953b605ac8f2bdf6aeffa/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht"
$configUrl   = "$repoBase/.localization-config"

function Set-Row($ws, $row, $values) {
    foreach ($pair in $values) {
        $col = $pair[0]
        $val = $pair[1]
        $cell = $ws.Range("$col$row")
        $cell.Value2 = $val
        if ($col -eq "D") {
            # Column D ("Latest Handoff Datetime") always carries the
            # custom date/time number format in this workbook.
            $cell.NumberFormat = "yyyy-mm-dd HH:mm:ss"
        }
    }
}

function Rebuild-Hyperlinks($ws, $links) {
    # Clear every hyperlink currently on the sheet, then add them back
    # in the desired final order so relationship ids come out sequential.
    $ws.Cells.Hyperlinks.Delete()
    foreach ($link in $links) {
        $cellRef = $link[0]
        $address = $link[1]
        $display = $link[2]
        $ws.Hyperlinks.Add($ws.Range($cellRef), $address, "", "", $display)
        $ws.Range($cellRef).Font.Underline = $true
        $ws.Range($cellRef).Font.Color = 15570276
    }
}

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

Set-Row $ws1 2 @(@("A",$mdFile), @("B","Ready for handoff"), @("C","Ready for handoff"))
Set-Row $ws1 3 @(@("A",$img1), @("B","Ready for handoff"), @("C","Ready for handoff"))
Set-Row $ws1 4 @(@("A",$img2), @("B","Ready for handoff"), @("C","Ready for handoff"))
Set-Row $ws1 5 @(@("A",".localization-config"), @("B","Not to be localized"), @("C","Not to be localized"))

Rebuild-Hyperlinks $ws1 @(
    @("A2", "$e2eBase/$mdFile", $mdFile),
    @("A3", "$e2eBase/$img1", $img1),
    @("A4", "$e2eBase/$img2", $img2),
    @("A5", $configUrl, ".localization-config")
)

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

Set-Row $ws2 2 @(@("A",$mdFile), @("B","Ready for handoff"), @("C",$xlfZh), @("D",$zhTime), @("G",$zeroDate), @("H","Include"))
Set-Row $ws2 3 @(@("A",$img1), @("B","Ready for handoff"), @("C",$img1Handoff), @("D",$zhTime), @("G",$zeroDate), @("H","IsDependency"), @("I",$depFrom))
Set-Row $ws2 4 @(@("A",$img2), @("B","Ready for handoff"), @("C",$img2Handoff), @("D",$zhTime), @("G",$zeroDate), @("H","IsDependency"), @("I",$depFrom))
Set-Row $ws2 5 @(@("A",".localization-config"), @("B","Not to be localized"), @("D",$zeroDate), @("G",$zeroDate), @("H","Ignored"))

Rebuild-Hyperlinks $ws2 @(
    @("A2", "$e2eBase/$mdFile", $mdFile),
    @("C2", "$handoffZhBase/$xlfZh", $xlfZh),
    @("A3", "$e2eBase/$img1", $img1),
    @("C3", "$handoffZhBase/$img1Handoff", $img1Handoff),
    @("A4", "$e2eBase/$img2", $img2),
    @("C4", "$handoffZhBase/$img2Handoff", $img2Handoff),
    @("A5", $configUrl, ".localization-config")
)

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)

Set-Row $ws3 2 @(@("A",$mdFile), @("B","Ready for handoff"), @("C",$xlfDe), @("D",$deTime), @("G",$zeroDate), @("H","Include"))
Set-Row $ws3 3 @(@("A",$img1), @("B","Ready for handoff"), @("C",$img1Handoff), @("D",$deTime), @("G",$zeroDate), @("H","IsDependency"), @("I",$depFrom))
Set-Row $ws3 4 @(@("A",$img2), @("B","Ready for handoff"), @("C",$img2Handoff), @("D",$deTime), @("G",$zeroDate), @("H","IsDependency"), @("I",$depFrom))
Set-Row $ws3 5 @(@("A",".localization-config"), @("B","Not to be localized"), @("D",$zeroDate), @("G",$zeroDate), @("H","Ignored"))

Rebuild-Hyperlinks $ws3 @(
    @("A2", "$e2eBase/$mdFile", $mdFile),
    @("C2", "$handoffDeBase/$xlfDe", $xlfDe),
    @("A3", "$e2eBase/$img1", $img1),
    @("C3", "$handoffDeBase/$img1Handoff", $img1Handoff),
    @("A4", "$e2eBase/$img2", $img2),
    @("C4", "$handoffDeBase/$img2Handoff", $img2Handoff),
    @("A5", $configUrl, ".localization-config")
)
